$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.615.18"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.596.66"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'211.12"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "1.820.22"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.615.62"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "'64.88"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "26.602.49"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'208.18"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "'7.00"
$ws.Range("E21").Value = "  +4.99%  "
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'2.31"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'145.43"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'7.10"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "'15.28"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("D34").Value = "1.284.06"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "'0.619"
$ws.Range("E35").Value = "  -6.36%  "
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").Value = "'0.840"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "'1.04"
$ws.Range("E40").Value = "  +19.30%  "
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'64.27"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.784"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("D45").Value = "1.732.79"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "'90.11"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("E48").Value = "  +4.28%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  -0.64%  "
